$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix row 245 C:F values (1568100000000 -> 1568092000000)
$ws.Range("C245:F245").Value = 1568092000000

# 2. Append new rows 253-255, reusing row 252's formatting
$ws.Range("A252:G252").Copy($ws.Range("A253:A253"))
$ws.Range("A252:G252").Copy($ws.Range("A254:A254"))
$ws.Range("A252:G252").Copy($ws.Range("A255:A255"))

$newRows = @(
    @{ Row = 253; Date = 44927.45833333334; C = 1719863000000 },
    @{ Row = 254; Date = 44958.45833333334; C = 1749400000000 },
    @{ Row = 255; Date = 44986.45833333334; C = 1788400000000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:AEM2"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.C
    $ws.Cells.Item($row, 5).Value = $r.C
    $ws.Cells.Item($row, 6).Value = $r.C
    $ws.Cells.Item($row, 7).Value = 0
}
